$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.9622023807194182
$arr[0,1] = 0.2218852548641905
$arr[0,2] = 0
$arr[0,3] = 0.2247505014119131
$arr[0,4] = 1.81982215889326
$arr[0,5] = 0.002443427777046538
$arr[0,6] = 0
$arr[0,7] = 0.7871052767708591
$arr[0,8] = 0.04259152284857493
$arr[0,9] = 0
$arr[0,10] = 0.4801065866548555
$arr[0,11] = 0.309811780404182
$arr[0,12] = 1.438640476772846
$arr[0,13] = 2.955193148474365
$ws.Range("B2:O2").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.8845360309094588
$arr[0,1] = 0.2120804661247178
$arr[0,2] = 0
$arr[0,3] = 0.2257689880360108
$arr[0,4] = 1.820648947197
$arr[0,5] = 0.002446034467612376
$arr[0,6] = 0
$arr[0,7] = 0.7938492463015372
$arr[0,8] = 0.04107525212612018
$arr[0,9] = 0
$arr[0,10] = 0.4753160814297388
$arr[0,11] = 0.2956323305692266
$arr[0,12] = 1.450440716834422
$arr[0,13] = 2.962596872369147
$ws.Range("B3:O3").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.836999538414176
$arr[0,1] = 0.2060009276992361
$arr[0,2] = 0
$arr[0,3] = 0.2264504418333395
$arr[0,4] = 1.82208039056858
$arr[0,5] = 0.00244772168579236
$arr[0,6] = 0
$arr[0,7] = 0.7984042883416755
$arr[0,8] = 0.04013807993391438
$arr[0,9] = 0
$arr[0,10] = 0.4725428870850408
$arr[0,11] = 0.2870185274497388
$arr[0,12] = 1.458192263878544
$arr[0,13] = 2.96890091769589
$ws.Range("B4:O4").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.817667351410762
$arr[0,1] = 0.2035086408486819
$arr[0,2] = 0
$arr[0,3] = 0.2267422766603424
$arr[0,4] = 1.822896201415119
$arr[0,5] = 0.00244843110657457
$arr[0,6] = 0
$arr[0,7] = 0.8003646437935465
$arr[0,8] = 0.03975464569393239
$arr[0,9] = 0
$arr[0,10] = 0.4714552622817507
$arr[0,11] = 0.2835318601884751
$arr[0,12] = 1.461478468903472
$arr[0,13] = 2.971911850174365
$ws.Range("B5:O5").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.8144596713469809
$arr[0,1] = 0.2030939060549457
$arr[0,2] = 0
$arr[0,3] = 0.2267915903360134
$arr[0,4] = 1.823045713398429
$arr[0,5] = 0.002448550227946594
$arr[0,6] = 0
$arr[0,7] = 0.8006964484878161
$arr[0,8] = 0.03969088509064989
$arr[0,9] = 0
$arr[0,10] = 0.4712772336687863
$arr[0,11] = 0.2829543310963913
$arr[0,12] = 1.462031837907546
$arr[0,13] = 2.972438504626496
$ws.Range("B6:O6").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.8367386567732922
$arr[0,1] = 0.2059673757406699
$arr[0,2] = 0
$arr[0,3] = 0.2264543203396396
$arr[0,4] = 1.822090451298855
$arr[0,5] = 0.002447731164665538
$arr[0,6] = 0
$arr[0,7] = 0.7984303047532713
$arr[0,8] = 0.04013291496422156
$arr[0,9] = 0
$arr[0,10] = 0.4725280467797432
$arr[0,11] = 0.2869714093979212
$arr[0,12] = 1.458236066856522
$arr[0,13] = 2.968939734766025
$ws.Range("B7:O7").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.9353924817924906
$arr[0,1] = 0.2185169411769863
$arr[0,2] = 0
$arr[0,3] = 0.225090050459837
$arr[0,4] = 1.819915612627412
$arr[0,5] = 0.002444308608153875
$arr[0,6] = 0
$arr[0,7] = 0.7893446147464971
$arr[0,8] = 0.04207000847975095
$arr[0,9] = 0
$arr[0,10] = 0.4784200173217954
$arr[0,11] = 0.3049036861023637
$arr[0,12] = 1.442604213729446
$arr[0,13] = 2.957380978266826
$ws.Range("B8:O8").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.130000184549658
$arr[0,1] = 0.2426522564813354
$arr[0,2] = 0
$arr[0,3] = 0.2228585144048125
$arr[0,4] = 1.822974334918783
$arr[0,5] = 0.002438281967583629
$arr[0,6] = 0
$arr[0,7] = 0.7748152880512791
$arr[0,8] = 0.04581875039774985
$arr[0,9] = 0
$arr[0,10] = 0.4913017484657161
$arr[0,11] = 0.3407917395345663
$arr[0,12] = 1.415961383288817
$arr[0,13] = 2.948672655134828
$ws.Range("B9:O9").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.27362549407934
$arr[0,1] = 0.260092343007301
$arr[0,2] = 0
$arr[0,3] = 0.221487779424379
$arr[0,4] = 1.829680853674219
$arr[0,5] = 0.002434267659509192
$arr[0,6] = 0
$arr[0,7] = 0.7661465637990617
$arr[0,8] = 0.04854158111374574
$arr[0,9] = 0
$arr[0,10] = 0.5015673617268135
$arr[0,11] = 0.3675880932563516
$arr[0,12] = 1.398825566617639
$arr[0,13] = 2.950800437349898
$ws.Range("B10:O10").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.339094564690186
$arr[0,1] = 0.2679623238522311
$arr[0,2] = 0
$arr[0,3] = 0.22092218560757
$arr[0,4] = 1.833699180361464
$arr[0,5] = 0.002432530351295894
$arr[0,6] = 0
$arr[0,7] = 0.762638913911843
$arr[0,8] = 0.04977326637811785
$arr[0,9] = 0
$arr[0,10] = 0.5064098801279755
$arr[0,11] = 0.3798694654691417
$arr[0,12] = 1.391558242272261
$arr[0,13] = 2.953623563164058
$ws.Range("B11:O11").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.363903962910058
$arr[0,1] = 0.2709332532007807
$arr[0,2] = 0
$arr[0,3] = 0.2207163146182172
$arr[0,4] = 1.835359804742254
$arr[0,5] = 0.002431885183926688
$arr[0,6] = 0
$arr[0,7] = 0.7613733592627199
$arr[0,8] = 0.05023865324113785
$arr[0,9] = 0
$arr[0,10] = 0.5082682783025376
$arr[0,11] = 0.3845330203646142
$arr[0,12] = 1.388882111583769
$arr[0,13] = 2.954959610295617
$ws.Range("B12:O12").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.358560052099449
$arr[0,1] = 0.2702938241618824
$arr[0,2] = 0
$arr[0,3] = 0.2207602835952009
$arr[0,4] = 1.83499598124321
$arr[0,5] = 0.002432023567655491
$arr[0,6] = 0
$arr[0,7] = 0.7616431289942511
$arr[0,8] = 0.05013846989922399
$arr[0,9] = 0
$arr[0,10] = 0.5078669459903153
$arr[0,11] = 0.383528072985527
$arr[0,12] = 1.389455091730028
$arr[0,13] = 2.954659990721694
$ws.Range("B13:O13").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.341135302859527
$arr[0,1] = 0.2682069302754257
$arr[0,2] = 0
$arr[0,3] = 0.2209050821476417
$arr[0,4] = 1.833833016621483
$arr[0,5] = 0.002432477018576152
$arr[0,6] = 0
$arr[0,7] = 0.7625335387602732
$arr[0,8] = 0.0498115747455472
$arr[0,9] = 0
$arr[0,10] = 0.5065622788630435
$arr[0,11] = 0.380252882821182
$arr[0,12] = 1.391336555604006
$arr[0,13] = 2.953728128438598
$ws.Range("B14:O14").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.330464399907839
$arr[0,1] = 0.2669274379482829
$arr[0,2] = 0
$arr[0,3] = 0.2209948564586544
$arr[0,4] = 1.833138760793702
$arr[0,5] = 0.002432756424201451
$arr[0,6] = 0
$arr[0,7] = 0.7630871094192386
$arr[0,8] = 0.04961120765554483
$arr[0,9] = 0
$arr[0,10] = 0.5057663361812672
$arr[0,11] = 0.37824840079859
$arr[0,12] = 1.392498882516826
$arr[0,13] = 2.953192112358806
$ws.Range("B15:O15").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.269349490771106
$arr[0,1] = 0.2595767315208093
$arr[0,2] = 0
$arr[0,3] = 0.2215259062666171
$arr[0,4] = 1.829437701984375
$arr[0,5] = 0.00243438297908953
$arr[0,6] = 0
$arr[0,7] = 0.7663845705150578
$arr[0,8] = 0.04846094554000757
$arr[0,9] = 0
$arr[0,10] = 0.501254349994511
$arr[0,11] = 0.3667872935602787
$arr[0,12] = 1.399311120050484
$arr[0,13] = 2.95065328913509
$ws.Range("B16:O16").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.23189054077136
$arr[0,1] = 0.2550509388368312
$arr[0,2] = 0
$arr[0,3] = 0.2218665135522322
$arr[0,4] = 1.827414912589092
$arr[0,5] = 0.002435403525002406
$arr[0,6] = 0
$arr[0,7] = 0.7685191090281833
$arr[0,8] = 0.04775349919953698
$arr[0,9] = 0
$arr[0,10] = 0.4985304789030778
$arr[0,11] = 0.3597795148803442
$arr[0,12] = 1.403625359599332
$arr[0,13] = 2.949571116180749
$ws.Range("B17:O17").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.210357768429503
$arr[0,1] = 0.2524418429730133
$arr[0,2] = 0
$arr[0,3] = 0.2220678790002122
$arr[0,4] = 1.826342532337989
$arr[0,5] = 0.002435998879737189
$arr[0,6] = 0
$arr[0,7] = 0.769787855382706
$arr[0,8] = 0.04734594278627924
$arr[0,9] = 0
$arr[0,10] = 0.4969800409174638
$arr[0,11] = 0.3557574636274268
$arr[0,12] = 1.406156482402594
$arr[0,13] = 2.94912328501951
$ws.Range("B18:O18").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.203069354919933
$arr[0,1] = 0.2515574244676486
$arr[0,2] = 0
$arr[0,3] = 0.2221369959912831
$arr[0,4] = 1.825995089649552
$arr[0,5] = 0.002436201894996057
$arr[0,6] = 0
$arr[0,7] = 0.7702244741805373
$arr[0,8] = 0.04720784003153256
$arr[0,9] = 0
$arr[0,10] = 0.4964578876348469
$arr[0,11] = 0.3543971589461492
$arr[0,12] = 1.407022012544509
$arr[0,13] = 2.949001639216078
$ws.Range("B19:O19").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.235876813976461
$arr[0,1] = 0.2555333368995321
$arr[0,2] = 0
$arr[0,3] = 0.2218296907311021
$arr[0,4] = 1.82762081719261
$arr[0,5] = 0.002435294020974643
$arr[0,6] = 0
$arr[0,7] = 0.7682876383727972
$arr[0,8] = 0.04782887574441474
$arr[0,9] = 0
$arr[0,10] = 0.4988187580502625
$arr[0,11] = 0.3605246131094688
$arr[0,12] = 1.403160959504227
$arr[0,13] = 2.949668242681639
$ws.Range("B20:O20").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.346252905924928
$arr[0,1] = 0.2688201538755379
$arr[0,2] = 0
$arr[0,3] = 0.2208623260688611
$arr[0,4] = 1.834170837111159
$arr[0,5] = 0.002432343484488726
$arr[0,6] = 0
$arr[0,7] = 0.7622703013349579
$arr[0,8] = 0.04990761982715952
$arr[0,9] = 0
$arr[0,10] = 0.5069448239234617
$arr[0,11] = 0.3812145388024533
$arr[0,12] = 1.390781866009419
$arr[0,13] = 2.953994591208954
$ws.Range("B21:O21").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.418492445733875
$arr[0,1] = 0.2774498196122295
$arr[0,2] = 0
$arr[0,3] = 0.2202785048707518
$arr[0,4] = 1.839261612666917
$arr[0,5] = 0.002430489216119866
$arr[0,6] = 0
$arr[0,7] = 0.7587032048056841
$arr[0,8] = 0.0512602085849494
$arr[0,9] = 0
$arr[0,10] = 0.5123992223483071
$arr[0,11] = 0.3948114525845696
$arr[0,12] = 1.383133485835884
$arr[0,13] = 2.95837843542688
$ws.Range("B22:O22").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.379927961389456
$arr[0,1] = 0.2728489871946636
$arr[0,2] = 0
$arr[0,3] = 0.220585681001511
$arr[0,4] = 1.836470509698216
$arr[0,5] = 0.002431472115188512
$arr[0,6] = 0
$arr[0,7] = 0.7605735643915423
$arr[0,8] = 0.05053886361051241
$arr[0,9] = 0
$arr[0,10] = 0.5094750336066625
$arr[0,11] = 0.3875477748791667
$arr[0,12] = 1.387175138057245
$arr[0,13] = 2.955896220134804
$ws.Range("B23:O23").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.234074610786138
$arr[0,1] = 0.2553152670223255
$arr[0,2] = 0
$arr[0,3] = 0.2218463210523165
$arr[0,4] = 1.827527445658575
$arr[0,5] = 0.002435343500846731
$arr[0,6] = 0
$arr[0,7] = 0.7683921567655148
$arr[0,8] = 0.04779480061066721
$arr[0,9] = 0
$arr[0,10] = 0.4986883785926892
$arr[0,11] = 0.3601877329520349
$arr[0,12] = 1.403370756514178
$arr[0,13] = 2.949623788793446
$ws.Range("B24:O24").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.077236219352073
$arr[0,1] = 0.2361741023894552
$arr[0,2] = 0
$arr[0,3] = 0.2234148767391062
$arr[0,4] = 1.821363383362566
$arr[0,5] = 0.002439839429095472
$arr[0,6] = 0
$arr[0,7] = 0.7783937102772818
$arr[0,8] = 0.04481005618989187
$arr[0,9] = 0
$arr[0,10] = 0.4876755913296051
$arr[0,11] = 0.3310068475526577
$arr[0,12] = 1.422740254072288
$arr[0,13] = 2.949532450298562
$ws.Range("B25:O25").Value = $arr

Write-Output "Done"